# Actualización 10 de Mayo
# Adds newly identified "Rescatables" (make-up exam) students to the
# "Rescatables" sheet, keeping the existing rows sorted alphabetically by
# Paterno (father's last name) for the 6AEM group, and appending three more
# students that belong to the 6BEV group at the bottom of the list.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

# --- Insert two new rows (CERVANTES, CRISTOBAL) right after CERON (row 6) ---
$ws.Rows.Item(7).Insert()
$ws.Rows.Item(7).Insert()

$ws.Cells.Item(7, 1).Value = 18330051920149
$ws.Cells.Item(7, 2).Value = "CERVANTES"
$ws.Cells.Item(7, 3).Value = "TENTZOHUA"
$ws.Cells.Item(7, 4).Value = "GONZALO"
$ws.Cells.Item(7, 5).Value = "TEMAS DE FILOSOFÍA"
$ws.Cells.Item(7, 6).Value = "6AEM"
$ws.Cells.Item(7, 7).Value = 2

$ws.Cells.Item(8, 1).Value = 18330051920152
$ws.Cells.Item(8, 2).Value = "CRISTOBAL"
$ws.Cells.Item(8, 3).Value = "ROMERO"
$ws.Cells.Item(8, 4).Value = "EDGAR ARMANDO"
$ws.Cells.Item(8, 5).Value = "TEMAS DE FILOSOFÍA"
$ws.Cells.Item(8, 6).Value = "6AEM"
$ws.Cells.Item(8, 7).Value = 2

# --- Insert one new row (MARTINEZ RAMIREZ) right after the existing
#     MARTINEZ ORTEGA row (now shifted down to row 12) ---
$ws.Rows.Item(13).Insert()

$ws.Cells.Item(13, 1).Value = 18330051920161
$ws.Cells.Item(13, 2).Value = "MARTINEZ"
$ws.Cells.Item(13, 3).Value = "RAMIREZ"
$ws.Cells.Item(13, 4).Value = "ALONSO ELIAS"
$ws.Cells.Item(13, 5).Value = "TEMAS DE FILOSOFÍA"
$ws.Cells.Item(13, 6).Value = "6AEM"
$ws.Cells.Item(13, 7).Value = 2

# --- Insert one new row (SANCHEZ TRUJILLO) right after RAMIREZ DOMINGUEZ
#     (now shifted down to row 15) ---
$ws.Rows.Item(16).Insert()

$ws.Cells.Item(16, 1).Value = 18330051920176
$ws.Cells.Item(16, 2).Value = "SANCHEZ"
$ws.Cells.Item(16, 3).Value = "TRUJILLO"
$ws.Cells.Item(16, 4).Value = "ERIK JAIR"
$ws.Cells.Item(16, 5).Value = "TEMAS DE FILOSOFÍA"
$ws.Cells.Item(16, 6).Value = "6AEM"
$ws.Cells.Item(16, 7).Value = 2

# --- Append three more rescatables belonging to the 6BEV group at the end
#     of the list (rows 21-23) ---
$ws.Cells.Item(21, 1).Value = 18330051920045
$ws.Cells.Item(21, 2).Value = "APARICIO"
$ws.Cells.Item(21, 3).Value = "ZUÑIGA"
$ws.Cells.Item(21, 4).Value = "JEAN PIERRE"
$ws.Cells.Item(21, 5).Value = "TEMAS DE FILOSOFÍA"
$ws.Cells.Item(21, 6).Value = "6BEV"
$ws.Cells.Item(21, 7).Value = 2

$ws.Cells.Item(22, 1).Value = 18330051920054
$ws.Cells.Item(22, 2).Value = "DE LA LUZ"
$ws.Cells.Item(22, 3).Value = "VELAZCO"
$ws.Cells.Item(22, 4).Value = "EMMANUEL"
$ws.Cells.Item(22, 5).Value = "TEMAS DE FILOSOFÍA"
$ws.Cells.Item(22, 6).Value = "6BEV"
$ws.Cells.Item(22, 7).Value = 2

$ws.Cells.Item(23, 1).Value = 18330051920080
$ws.Cells.Item(23, 2).Value = "SANCHEZ"
$ws.Cells.Item(23, 3).Value = "ALMANZA"
$ws.Cells.Item(23, 4).Value = "MARIO"
$ws.Cells.Item(23, 5).Value = "TEMAS DE FILOSOFÍA"
$ws.Cells.Item(23, 6).Value = "6BEV"
$ws.Cells.Item(23, 7).Value = 2
